$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matches source formatting)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated coin data
$ws.Range("D2").Value = "43.529.29"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "2.245.33"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "230.50"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").Value = "0.641"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("D7").Value = "63.68"
$ws.Range("E7").Value = "  -3.28%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("D10").Value = "0.0958"
$ws.Range("E10").Value = "  -7.90%  "
$ws.Range("D11").Value = "56.89"
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").Value = "27.73"
$ws.Range("E12").Value = "  +5.60%  "
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "2.578.19"
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").Value = "15.45"
$ws.Range("E15").Value = "  -3.14%  "
$ws.Range("D16").Value = "6.08"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "2.240.02"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("D19").Value = "43.432.21"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").Value = "0.0₃0966"
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("D21").Value = "72.87"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("D23").Value = "246.01"
$ws.Range("E23").Value = "  -5.98%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "3.69"
$ws.Range("E25").Value = "  +31.39%  "
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "9.75"
$ws.Range("E28").Value = "  -5.09%  "
$ws.Range("D29").Value = "172.68"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "21.50"
$ws.Range("E30").Value = "  +2.01%  "
$ws.Range("E31").Value = "  -4.96%  "
$ws.Range("D32").Value = "1.42"
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("E34").Value = "  +2.54%  "
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("D36").Value = "4.87"
$ws.Range("E36").Value = "  -3.56%  "
$ws.Range("D37").Value = "3.61"
$ws.Range("E37").Value = "  -6.73%  "
$ws.Range("E38").Value = "  -8.34%  "
$ws.Range("E39").Value = "  -4.46%  "
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "8.62"
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("D43").Value = "4.50"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "17.00"
$ws.Range("E44").Value = "  -3.98%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "96.38"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.0938"
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("B47").Value = "TerraClassic"
$ws.Range("C47").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D47").Value = "0.000209"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("E48").Value = "  -2.38%  "
$ws.Range("D49").Value = "1.449.53"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("D50").Value = "9.98"
$ws.Range("E50").Value = "  -3.50%  "
$ws.Range("D51").Value = "2.29"
$ws.Range("E51").Value = "  -3.02%  "
